$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GlobalConstantIntTable")

# ---------------------------------------------------------------
# 1. Insert a new row at position 30 (pushes old rows 30-36 down to 31-37)
# ---------------------------------------------------------------
$ws.Rows.Item(30).Insert()

# ---------------------------------------------------------------
# 2. Append two brand-new rows at the bottom of the table
# ---------------------------------------------------------------
$ws.Range("A38").Value() = "PetSaleGivenTime"
$ws.Range("A39").Value() = "PetSaleCoolTime"

# Row 31 (was old row 30: PetDailySearchCount) now also becomes a
# "hard-coded" style row with a 하드코딩 note in column B
$ws.Range("B31").Value() = "하드코딩"
$ws.Range("A31").Font.Color = $ws.Range("A9").Font.Color

# Row 36 (was old row 35: PetHeartCount) gets the same treatment
$ws.Range("B36").Value() = "하드코딩"
$ws.Range("A36").Font.Color = $ws.Range("A9").Font.Color

$ws.Range("C38").Value() = "펫 세일 열리는 시간"
$ws.Range("C39").Value() = "펫 세일 쿨타임"
$ws.Range("D38").Value() = 86400
$ws.Range("D39").Value() = 172800

# New row 30: MaxPetCountStep
$ws.Range("C30").Value() = "구매 가능 최대마리수 스텝"
$ws.Range("A30").Value() = "MaxPetCountStep"
$ws.Range("B30").Value() = "int"
$ws.Range("D30").Value() = 9
# match the "hard-coded variable" styling used by similar rows (blue font)
$ws.Range("A30").Font.Color = $ws.Range("A9").Font.Color

# ---------------------------------------------------------------
# 3. Extend the big JSON-building formula in F2 to also emit MaxPetCountStep
# ---------------------------------------------------------------
$newFormula = '="{"""&A2&""":"&VLOOKUP(A2,$A:$D,MATCH($D$1,$A$1:$D$1,0),0)&","""&A9&""":"&VLOOKUP(A9,$A:$D,MATCH($D$1,$A$1:$D$1,0),0)&","""&A10&""":"&VLOOKUP(A10,$A:$D,MATCH($D$1,$A$1:$D$1,0),0)&","""&A11&""":"&VLOOKUP(A11,$A:$D,MATCH($D$1,$A$1:$D$1,0),0)&","""&A12&""":"&VLOOKUP(A12,$A:$D,MATCH($D$1,$A$1:$D$1,0),0)&","""&A13&""":"&VLOOKUP(A13,$A:$D,MATCH($D$1,$A$1:$D$1,0),0)&","""&A14&""":"&VLOOKUP(A14,$A:$D,MATCH($D$1,$A$1:$D$1,0),0)&","""&A15&""":"&VLOOKUP(A15,$A:$D,MATCH($D$1,$A$1:$D$1,0),0)&","""&A18&""":"&VLOOKUP(A18,$A:$D,MATCH($D$1,$A$1:$D$1,0),0)&","""&A27&""":"&VLOOKUP(A27,$A:$D,MATCH($D$1,$A$1:$D$1,0),0)&","""&A28&""":"&VLOOKUP(A28,$A:$D,MATCH($D$1,$A$1:$D$1,0),0)&","""&A29&""":"&VLOOKUP(A29,$A:$D,MATCH($D$1,$A$1:$D$1,0),0)&","""&A30&""":"&VLOOKUP(A30,$A:$D,MATCH($D$1,$A$1:$D$1,0),0)&"}"'
$ws.Range("F2").Formula = $newFormula

# ---------------------------------------------------------------
# 4. Update the view/selection state to match the authored edit
# ---------------------------------------------------------------
$ws.Activate()
$ws.Range("A34").Select()

Write-Host "Done applying pet constant edits"
